$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_1_7_24"
$ws.Range("B2").Value = 0.9811311972403775
$ws.Range("C2").Value = 0.9722766314271021
$ws.Range("D2").Value = 0.8607819149208679
$ws.Range("E2").Value = 0.9312733809837527
$ws.Range("F2").Value = 2.275583230630684
$ws.Range("G2").Value = 3.293514223747568
$ws.Range("H2").Value = 12.19482068099679
$ws.Range("I2").Value = 7.482360770893265

$ws.Range("A3").Value = "model_1_7_23"
$ws.Range("B3").Value = 0.9813298165793169
$ws.Range("C3").Value = 0.9724884863666445
$ws.Range("D3").Value = 0.8624507479474339
$ws.Range("E3").Value = 0.9320276139178394
$ws.Range("F3").Value = 2.251629679219512
$ws.Range("G3").Value = 3.268346024763396
$ws.Range("H3").Value = 12.04863910197327
$ws.Range("I3").Value = 7.400246402415569

$ws.Range("A4").Value = "model_1_7_22"
$ws.Range("B4").Value = 0.9815472214928532
$ws.Range("C4").Value = 0.9727224996571464
$ws.Range("D4").Value = 0.8642853424337547
$ws.Range("E4").Value = 0.9328573431207275
$ws.Range("F4").Value = 2.225410581918938
$ws.Range("G4").Value = 3.240545431239303
$ws.Range("H4").Value = 11.88793763297727
$ws.Range("I4").Value = 7.309912652159543

$ws.Range("A5").Value = "model_1_7_21"
$ws.Range("B5").Value = 0.9817848154590446
$ws.Range("C5").Value = 0.9729807355264543
$ws.Range("D5").Value = 0.8662996990820115
$ws.Range("E5").Value = 0.9337692619209125
$ws.Range("F5").Value = 2.19675667885727
$ws.Range("G5").Value = 3.209867214542399
$ws.Range("H5").Value = 11.71148988124225
$ws.Range("I5").Value = 7.210630808320685

$ws.Range("A6").Value = "model_1_7_20"
$ws.Range("B6").Value = 0.9820438491242501
$ws.Range("C6").Value = 0.9732655324207524
$ws.Range("D6").Value = 0.8685081299685756
$ws.Range("E6").Value = 0.9347699614602124
$ws.Range("F6").Value = 2.165517141711232
$ws.Range("G6").Value = 3.176033569118561
$ws.Range("H6").Value = 11.51804217915157
$ws.Range("I6").Value = 7.101683284297446

$ws.Range("A7").Value = "model_1_7_19"
$ws.Range("B7").Value = 0.9823255699958146
$ws.Range("C7").Value = 0.9735791442076972
$ws.Range("D7").Value = 0.8709257606109757
$ws.Range("E7").Value = 0.9358664194534637
$ws.Range("F7").Value = 2.131541520723631
$ws.Range("G7").Value = 3.138776737275712
$ws.Range("H7").Value = 11.30626960563719
$ws.Range("I7").Value = 6.982310406756407

$ws.Range("A8").Value = "model_1_7_18"
$ws.Range("B8").Value = 0.9826311098352474
$ws.Range("C8").Value = 0.973924504405443
$ws.Range("D8").Value = 0.8735666033776699
$ws.Range("E8").Value = 0.9370659280770955
$ws.Range("F8").Value = 2.094693325119465
$ws.Range("G8").Value = 3.097748219381111
$ws.Range("H8").Value = 11.07494474602402
$ws.Range("I8").Value = 6.851718266501568

$ws.Range("A9").Value = "model_1_7_17"
$ws.Range("B9").Value = 0.9829612907920826
$ws.Range("C9").Value = 0.9743039581453855
$ws.Range("D9").Value = 0.8764457559725447
$ws.Range("E9").Value = 0.9383751746373016
$ws.Range("F9").Value = 2.054873403420164
$ws.Range("G9").Value = 3.052669415682741
$ws.Range("H9").Value = 10.82274511558257
$ws.Range("I9").Value = 6.709178807384616

$ws.Range("A10").Value = "model_1_7_16"
$ws.Range("B10").Value = 0.9833167874654742
$ws.Range("C10").Value = 0.9747204811282001
$ws.Range("D10").Value = 0.8795767010833151
$ws.Range("E10").Value = 0.9398012853718554
$ws.Range("F10").Value = 2.012000398766892
$ws.Range("G10").Value = 3.003186815297787
$ws.Range("H10").Value = 10.54848969707005
$ws.Range("I10").Value = 6.553916186177356

$ws.Range("A11").Value = "model_1_7_15"
$ws.Range("B11").Value = 0.9836977117555801
$ws.Range("C11").Value = 0.9751768905462379
$ws.Range("D11").Value = 0.8829718144012375
$ws.Range("E11").Value = 0.9413504381097698
$ws.Range("F11").Value = 1.966060815967308
$ws.Range("G11").Value = 2.948965738006723
$ws.Range("H11").Value = 10.25109444069804
$ws.Range("I11").Value = 6.38525781420721

$ws.Range("A12").Value = "model_1_7_14"
$ws.Range("B12").Value = 0.9841036422529196
$ws.Range("C12").Value = 0.9756760605033875
$ws.Range("D12").Value = 0.8866415272025865
$ws.Range("E12").Value = 0.9430282031591959
$ws.Range("F12").Value = 1.917105477130225
$ws.Range("G12").Value = 2.889664742544463
$ws.Range("H12").Value = 9.929645617883304
$ws.Range("I12").Value = 6.202597244426614

$ws.Range("A13").Value = "model_1_7_13"
$ws.Range("B13").Value = 0.9845333848089411
$ws.Range("C13").Value = 0.9762207932395635
$ws.Range("D13").Value = 0.8905929237580978
$ws.Range("E13").Value = 0.9448389216431242
$ws.Range("F13").Value = 1.865278396926508
$ws.Range("G13").Value = 2.824950925029157
$ws.Range("H13").Value = 9.583522681293717
$ws.Range("I13").Value = 6.005461852853324

$ws.Range("A14").Value = "model_1_7_12"
$ws.Range("B14").Value = 0.9849846514956928
$ws.Range("C14").Value = 0.9768137069052711
$ws.Range("D14").Value = 0.8948283689205625
$ws.Range("E14").Value = 0.9467851923299563
$ws.Range("F14").Value = 1.810855500148348
$ws.Range("G14").Value = 2.754513251254845
$ws.Range("H14").Value = 9.212518481436387
$ws.Range("I14").Value = 5.793568708026148

$ws.Range("A15").Value = "model_1_7_11"
$ws.Range("B15").Value = 0.985453819833088
$ws.Range("C15").Value = 0.9774570264991724
$ws.Range("D15").Value = 0.8993454370757346
$ws.Range("E15").Value = 0.9488670366475529
$ws.Range("F15").Value = 1.754273658972721
$ws.Range("G15").Value = 2.678087393143191
$ws.Range("H15").Value = 8.816845490209328
$ws.Range("I15").Value = 5.566915477064649

$ws.Range("A16").Value = "model_1_7_10"
$ws.Range("B16").Value = 0.9859355023992954
$ws.Range("C16").Value = 0.9781523294398173
$ws.Range("D16").Value = 0.9041328901727765
$ws.Range("E16").Value = 0.9510813514548569
$ws.Range("F16").Value = 1.696182598076463
$ws.Range("G16").Value = 2.595485954620088
$ws.Range("H16").Value = 8.397488105686156
$ws.Range("I16").Value = 5.325840003169116

$ws.Range("A17").Value = "model_1_7_9"
$ws.Range("B17").Value = 0.9864220142641329
$ws.Range("C17").Value = 0.9789006127597101
$ws.Range("D17").Value = 0.9091692465086308
$ws.Range("E17").Value = 0.9534205015585657
$ws.Range("F17").Value = 1.637509122327581
$ws.Range("G17").Value = 2.506590489013901
$ws.Range("H17").Value = 7.956328019577834
$ws.Range("I17").Value = 5.071173540250526

$ws.Range("A18").Value = "model_1_7_8"
$ws.Range("B18").Value = 0.9869028579635222
$ws.Range("C18").Value = 0.9797013445863284
$ws.Range("D18").Value = 0.9144204055068657
$ws.Range("E18").Value = 0.9558713214810872
$ws.Range("F18").Value = 1.579519229019347
$ws.Range("G18").Value = 2.411464182358904
$ws.Range("H18").Value = 7.496352274943242
$ws.Range("I18").Value = 4.804349431814988

$ws.Range("A19").Value = "model_1_7_7"
$ws.Range("B19").Value = 0.987363886676932
$ws.Range("C19").Value = 0.9805521721951526
$ws.Range("D19").Value = 0.9198360323975835
$ws.Range("E19").Value = 0.9584132722547085
$ws.Range("F19").Value = 1.523919028919753
$ws.Range("G19").Value = 2.310386536464219
$ws.Range("H19").Value = 7.02196995047764
$ws.Range("I19").Value = 4.527603783297207

$ws.Range("A20").Value = "model_1_7_0"
$ws.Range("B20").Value = 0.987553816147049
$ws.Range("C20").Value = 0.9868362268560323
$ws.Range("D20").Value = 0.9523925687787472
$ws.Range("E20").Value = 0.974370185148775
$ws.Range("F20").Value = 1.501013478275832
$ws.Range("G20").Value = 1.56384582103876
$ws.Range("H20").Value = 4.170177218685851
$ws.Range("I20").Value = 2.790352907695453

$ws.Range("A21").Value = "model_1_7_6"
$ws.Range("B21").Value = 0.9877864504290333
$ws.Range("C21").Value = 0.9814484266773384
$ws.Range("D21").Value = 0.9253445672437212
$ws.Range("E21").Value = 0.9610166725174437
$ws.Range("F21").Value = 1.47295771460617
$ws.Range("G21").Value = 2.203912214001745
$ws.Range("H21").Value = 6.539449345302753
$ws.Range("I21").Value = 4.244168045068657

$ws.Range("A22").Value = "model_1_7_5"
$ws.Range("B22").Value = 0.9881461820351606
$ws.Range("C22").Value = 0.9823815268915742
$ws.Range("D22").Value = 0.9308493476109969
$ws.Range("E22").Value = 0.9636400191400931
$ws.Range("F22").Value = 1.429573975804111
$ws.Range("G22").Value = 2.093060647761275
$ws.Range("H22").Value = 6.057257614041385
$ws.Range("I22").Value = 3.958560719424897

$ws.Range("A23").Value = "model_1_7_1"
$ws.Range("B23").Value = 0.9881860326579773
$ws.Range("C23").Value = 0.9860982485051583
$ws.Range("D23").Value = 0.9496774066354537
$ws.Range("E23").Value = 0.9729158351899441
$ws.Range("F23").Value = 1.424767978827666
$ws.Range("G23").Value = 1.651517064489224
$ws.Range("H23").Value = 4.408012090775021
$ws.Range("I23").Value = 2.948689971774432

$ws.Range("A24").Value = "model_1_7_4"
$ws.Range("B24").Value = 0.9884116862444768
$ws.Range("C24").Value = 0.983338242428653
$ws.Range("D24").Value = 0.936222823145347
$ws.Range("E24").Value = 0.9662271770107898
$ws.Range("F24").Value = 1.39755408911184
$ws.Range("G24").Value = 1.979403599875341
$ws.Range("H24").Value = 5.58656754142707
$ws.Range("I24").Value = 3.676893312575857

$ws.Range("A25").Value = "model_1_7_2"
$ws.Range("B25").Value = 0.9884890070310878
$ws.Range("C25").Value = 0.9852324454321029
$ws.Range("D25").Value = 0.9458722818188694
$ws.Range("E25").Value = 0.9709749602194466
$ws.Range("F25").Value = 1.388229179225799
$ws.Range("G25").Value = 1.754373783670885
$ws.Range("H25").Value = 4.741322341240553
$ws.Range("I25").Value = 3.159995677603278

$ws.Range("A26").Value = "model_1_7_3"
$ws.Range("B26").Value = 0.9885428988183222
$ws.Range("C26").Value = 0.9842982911912282
$ws.Range("D26").Value = 0.9413002379540748
$ws.Range("E26").Value = 0.9687042249899793
$ws.Range("F26").Value = 1.381729813640083
$ws.Range("G26").Value = 1.865350567441038
$ws.Range("H26").Value = 5.141810934695422
$ws.Range("I26").Value = 3.407213719829898

